# I added 2 functions: delete_user and most_mistakes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DB")

$mistakes = @(1, 5, 2, 7, 4, 10, 0, 3, 14, 4, 8, 12, 9, 7, 6)

for ($i = 0; $i -lt $mistakes.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $mistakes[$i]
}

# Reflect the author's final cursor position in the sheet view
[void]$ws.Range("F10").Select()
